$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 23061.076
$ws.Range("I64").Value = 3808.348
$ws.Range("K64").Value = 3808.348
$ws.Range("M64").Value = -3560.348
$ws.Range("H67").Value = 23061.076
$ws.Range("I67").Value = 3808.348
$ws.Range("K67").Value = 3808.348
$ws.Range("M67").Value = -2950.348
$ws.Range("H81").Value = 43166.668
$ws.Range("J81").Value = 49750
$ws.Range("L81").Value = 49750
$ws.Range("N81").Value = -51746
$ws.Range("H84").Value = 43166.668
$ws.Range("J84").Value = 49750
$ws.Range("L84").Value = 149250
$ws.Range("N84").Value = -159234
$ws.Range("H100").Value = 6557.2144
$ws.Range("I100").Value = 1679.5
$ws.Range("K100").Value = 1679.5
$ws.Range("M100").Value = -1138.5
$ws.Range("H141").Value = 2500.1428
$ws.Range("I141").Value = 2069.5652
$ws.Range("K141").Value = 6208.6956
$ws.Range("M141").Value = -1028.6956

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2037
$ws.Range("I2").Value = 1211
$ws.Range("J2").Value = 2450
$ws.Range("K2").Value = 1211
$ws.Range("L2").Value = 2450
$ws.Range("M2").Value = -1098
$ws.Range("N2").Value = -2676
$ws.Range("H63").Value = 3222.2222
$ws.Range("I63").Value = 3000
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 3000
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -2314
$ws.Range("N63").Value = -6372
$ws.Range("H66").Value = 3222.2222
$ws.Range("I66").Value = 3000
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 15000
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -11568
$ws.Range("N66").Value = -31864
$ws.Range("H102").Value = 2497.5
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 2497.5
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 2497.5
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -5741.5
$ws.Range("H116").Value = 2037
$ws.Range("I116").Value = 1211
$ws.Range("J116").Value = 2450
$ws.Range("K116").Value = 1211
$ws.Range("L116").Value = 2450
$ws.Range("M116").Value = 1083
$ws.Range("N116").Value = -7038
$ws.Range("H122").Value = 1103.3684
$ws.Range("I122").Value = 943.38464
$ws.Range("K122").Value = 2830.15392
$ws.Range("M122").Value = -380.1539199999997

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2037
$ws.Range("I3").Value = 1211
$ws.Range("J3").Value = 2450
$ws.Range("K3").Value = 1211
$ws.Range("L3").Value = 2450
$ws.Range("M3").Value = -1097
$ws.Range("N3").Value = -2678
$ws.Range("H20").Value = 1323.5454
$ws.Range("I20").Value = 1273.2222
$ws.Range("J20").Value = 1550
$ws.Range("K20").Value = 1273.2222
$ws.Range("L20").Value = 1550
$ws.Range("M20").Value = -1026.2222
$ws.Range("N20").Value = -2044
$ws.Range("H99").Value = 1375
$ws.Range("I99").Value = 1200
$ws.Range("J99").Value = 1400
$ws.Range("K99").Value = 1200
$ws.Range("L99").Value = 1400
$ws.Range("M99").Value = 298
$ws.Range("N99").Value = -4396
$ws.Range("H105").Value = 443556.12
$ws.Range("I105").Value = 725047.25
$ws.Range("J105").Value = 1212.8572
$ws.Range("K105").Value = 725047.25
$ws.Range("L105").Value = 1212.8572
$ws.Range("M105").Value = -723300.25
$ws.Range("N105").Value = -4706.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22712.871
$ws.Range("I31").Value = 22258.893
$ws.Range("J31").Value = 26950
$ws.Range("K31").Value = 22258.893
$ws.Range("L31").Value = 26950
$ws.Range("M31").Value = -21963.893
$ws.Range("N31").Value = -27540
$ws.Range("H34").Value = 22712.871
$ws.Range("I34").Value = 22258.893
$ws.Range("J34").Value = 26950
$ws.Range("K34").Value = 22258.893
$ws.Range("L34").Value = 26950
$ws.Range("M34").Value = -22056.893
$ws.Range("N34").Value = -27354
$ws.Range("H99").Value = 2890.2693
$ws.Range("I99").Value = 2363.8333
$ws.Range("J99").Value = 4074.75
$ws.Range("K99").Value = 2363.8333
$ws.Range("L99").Value = 4074.75
$ws.Range("M99").Value = -865.8332999999998
$ws.Range("N99").Value = -7070.75
$ws.Range("H105").Value = 881.25
$ws.Range("I105").Value = 876.6667
$ws.Range("J105").Value = 895
$ws.Range("K105").Value = 876.6667
$ws.Range("L105").Value = 895
$ws.Range("M105").Value = 870.3333
$ws.Range("N105").Value = -4389
$ws.Range("H126").Value = 2890.2693
$ws.Range("I126").Value = 2363.8333
$ws.Range("J126").Value = 4074.75
$ws.Range("K126").Value = 7091.499899999999
$ws.Range("L126").Value = 12224.25
$ws.Range("M126").Value = -4621.499899999999
$ws.Range("N126").Value = -17164.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6194.8335
$ws.Range("I5").Value = 589.3
$ws.Range("J5").Value = 13201.75
$ws.Range("K5").Value = 1767.9
$ws.Range("L5").Value = 39605.25
$ws.Range("M5").Value = -1655.9
$ws.Range("N5").Value = -39829.25
$ws.Range("H122").Value = 499.33334
$ws.Range("I122").Value = 326.92856
$ws.Range("K122").Value = 2942.35704
$ws.Range("M122").Value = -492.3570399999999
$ws.Range("H132").Value = 861.7143
$ws.Range("I132").Value = 694.75
$ws.Range("J132").Value = 1084.3334
$ws.Range("K132").Value = 6252.75
$ws.Range("L132").Value = 9759.000599999999
$ws.Range("M132").Value = -3722.75
$ws.Range("N132").Value = -14819.0006
$ws.Range("H135").Value = 6194.8335
$ws.Range("I135").Value = 589.3
$ws.Range("J135").Value = 13201.75
$ws.Range("K135").Value = 5303.7
$ws.Range("L135").Value = 118815.75
$ws.Range("M135").Value = -2768.7
$ws.Range("N135").Value = -123885.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 58581.285
$ws.Range("J52").Value = 58581.285
$ws.Range("L52").Value = 58581.285
$ws.Range("N52").Value = -59099.285
$ws.Range("H70").Value = 4452.684
$ws.Range("I70").Value = 4314.357
$ws.Range("J70").Value = 4840
$ws.Range("K70").Value = 4314.357
$ws.Range("L70").Value = 4840
$ws.Range("M70").Value = -4044.357
$ws.Range("N70").Value = -5380
$ws.Range("H73").Value = 4452.684
$ws.Range("I73").Value = 4314.357
$ws.Range("J73").Value = 4840
$ws.Range("K73").Value = 4314.357
$ws.Range("L73").Value = 4840
$ws.Range("M73").Value = -3378.357
$ws.Range("N73").Value = -6712
